$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.852.89"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.588.26"
$ws.Range("E3").Value = "  -1.83%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.64"
$ws.Range("E5").Value = "  -1.38%  "

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.483"
$ws.Range("E7").Value = "  -3.31%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.247"
$ws.Range("E8").Value = "  -0.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0617"
$ws.Range("E9").Value = "  +0.09%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.14"
$ws.Range("E10").Value = "  -0.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.811.55"
$ws.Range("E12").Value = "  -1.63%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.588.59"
$ws.Range("E13").Value = "  -1.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  -2.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.511"
$ws.Range("E15").Value = "  -2.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.865.89"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0₃0723"
$ws.Range("E17").Value = "  -1.63%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.03"
$ws.Range("E18").Value = "  -2.22%  "

$ws.Range("E19").Value = "  -0.09%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.43"
$ws.Range("E20").Value = "  +0.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.19"
$ws.Range("E21").Value = "  -1.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.39"
$ws.Range("E22").Value = "  -1.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").Value = "  -1.55%  "

$ws.Range("E24").Value = "  -1.80%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.46"
$ws.Range("E25").Value = "  -1.50%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.14"
$ws.Range("E28").Value = "  -0.55%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.46"
$ws.Range("E29").Value = "  -2.89%  "

$ws.Range("E30").Value = "  -5.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0472"
$ws.Range("E31").Value = "  -1.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.12"
$ws.Range("E32").Value = "  -0.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.04"
$ws.Range("E33").Value = "  -1.72%  "

$ws.Range("E34").Value = "  +0.68%  "

$ws.Range("E35").Value = "  -2.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.103.64"
$ws.Range("E36").Value = "  -2.00%  "

$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.34"
$ws.Range("E38").Value = "  -1.68%  "

$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.505"
$ws.Range("E39").Value = "  -1.08%  "

$ws.Range("E40").Value = "  -1.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.777"
$ws.Range("E41").Value = "  -7.38%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.818"
$ws.Range("E42").Value = "  +9.17%  "

$ws.Range("E43").Value = "  +2.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "93.67"
$ws.Range("E44").Value = "  -4.50%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.724.58"
$ws.Range("E45").Value = "  -1.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₆0112"
$ws.Range("E46").Value = "  -0.40%  "

$ws.Range("E47").Value = "  +0.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.44"
$ws.Range("E48").Value = "  -1.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0509"
$ws.Range("E49").Value = "  -1.67%  "

$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("E51").Value = "  +0.00%  "

